$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing row (162) down to the new rows (163-178)
$ws.Range("A162:P162").Copy()
$ws.Range("A163:P178").PasteSpecial(-4122)
$ws.Range("163:178").RowHeight = 12.75

# FORNECEDOR_CDG (col O) values are zero-padded numeric codes; temporarily force
# text format so the leading zeros are preserved when the values are typed in below,
# then restore the original (General) cell format/style to match the rest of the column.
$ws.Range("O163:O178").NumberFormat = "@"

# Row 163
$ws.Cells.Item(163, 1).Value = 2510
$ws.Cells.Item(163, 2).Value = "SAMAUMA EVENTOS LTDA"
$ws.Cells.Item(163, 3).Value = "RJ"
$ws.Cells.Item(163, 4).Value = 24
$ws.Cells.Item(163, 5).Value = 46014.4998896181
$ws.Cells.Item(163, 6).Value = 81393
$ws.Cells.Item(163, 7).Value = 46014
$ws.Cells.Item(163, 8).Value = "E.01.0031"
$ws.Cells.Item(163, 9).Value = "MATERIAL DE ED AR CONDICIONADO"
$ws.Cells.Item(163, 10).Value = "Apto"
$ws.Cells.Item(163, 11).Value = "UN"
$ws.Cells.Item(163, 12).Value = 1
$ws.Cells.Item(163, 13).Value = 98031
$ws.Cells.Item(163, 14).Value = 98031
$ws.Cells.Item(163, 15).Value = "00000000006292"
$ws.Cells.Item(163, 16).Value = "DAIKIN"

# Row 164
$ws.Cells.Item(164, 1).Value = 2212
$ws.Cells.Item(164, 2).Value = "IDEA INVEST. IMOBILIÁRIOS LTDA."
$ws.Cells.Item(164, 3).Value = "RJ"
$ws.Cells.Item(164, 4).Value = 267
$ws.Cells.Item(164, 5).Value = 46014.6609215394
$ws.Cells.Item(164, 6).Value = 81399
$ws.Cells.Item(164, 7).Value = 46014
$ws.Cells.Item(164, 8).Value = "P2.02.0037"
$ws.Cells.Item(164, 9).Value = " MATERIAL PARA  ESQUADRIAS DE ALUMÍNIO ED -"
$ws.Cells.Item(164, 10).Value = "Apto"
$ws.Cells.Item(164, 11).Value = "VB"
$ws.Cells.Item(164, 12).Value = 1
$ws.Cells.Item(164, 13).Value = 88155
$ws.Cells.Item(164, 14).Value = 88155
$ws.Cells.Item(164, 15).Value = "00000000002083"
$ws.Cells.Item(164, 16).Value = "ALU-SERVIÇE"

# Row 165
$ws.Cells.Item(165, 1).Value = 2212
$ws.Cells.Item(165, 2).Value = "IDEA INVEST. IMOBILIÁRIOS LTDA."
$ws.Cells.Item(165, 3).Value = "RJ"
$ws.Cells.Item(165, 4).Value = 268
$ws.Cells.Item(165, 5).Value = 46014.6653912153
$ws.Cells.Item(165, 6).Value = 81402
$ws.Cells.Item(165, 7).Value = 46014
$ws.Cells.Item(165, 8).Value = "P2.02.0037"
$ws.Cells.Item(165, 9).Value = " MATERIAL PARA  ESQUADRIAS DE ALUMÍNIO ED -"
$ws.Cells.Item(165, 10).Value = "Apto"
$ws.Cells.Item(165, 11).Value = "VB"
$ws.Cells.Item(165, 12).Value = 1
$ws.Cells.Item(165, 13).Value = 87645
$ws.Cells.Item(165, 14).Value = 87645
$ws.Cells.Item(165, 15).Value = "00000000002083"
$ws.Cells.Item(165, 16).Value = "ALU-SERVIÇE"

# Row 166
$ws.Cells.Item(166, 1).Value = 2212
$ws.Cells.Item(166, 2).Value = "IDEA INVEST. IMOBILIÁRIOS LTDA."
$ws.Cells.Item(166, 3).Value = "RJ"
$ws.Cells.Item(166, 4).Value = 269
$ws.Cells.Item(166, 5).Value = 46014.6654302662
$ws.Cells.Item(166, 6).Value = 81401
$ws.Cells.Item(166, 7).Value = 46014
$ws.Cells.Item(166, 8).Value = "P2.02.0037"
$ws.Cells.Item(166, 9).Value = " MATERIAL PARA  ESQUADRIAS DE ALUMÍNIO ED -"
$ws.Cells.Item(166, 10).Value = "Apto"
$ws.Cells.Item(166, 11).Value = "VB"
$ws.Cells.Item(166, 12).Value = 1
$ws.Cells.Item(166, 13).Value = 58770
$ws.Cells.Item(166, 14).Value = 58770
$ws.Cells.Item(166, 15).Value = "00000000002083"
$ws.Cells.Item(166, 16).Value = "ALU-SERVIÇE"

# Row 167
$ws.Cells.Item(167, 1).Value = 2212
$ws.Cells.Item(167, 2).Value = "IDEA INVEST. IMOBILIÁRIOS LTDA."
$ws.Cells.Item(167, 3).Value = "RJ"
$ws.Cells.Item(167, 4).Value = 270
$ws.Cells.Item(167, 5).Value = 46014.6654644444
$ws.Cells.Item(167, 6).Value = 81403
$ws.Cells.Item(167, 7).Value = 46014
$ws.Cells.Item(167, 8).Value = "P2.02.0037"
$ws.Cells.Item(167, 9).Value = " MATERIAL PARA  ESQUADRIAS DE ALUMÍNIO ED -"
$ws.Cells.Item(167, 10).Value = "Apto"
$ws.Cells.Item(167, 11).Value = "VB"
$ws.Cells.Item(167, 12).Value = 1
$ws.Cells.Item(167, 13).Value = 40000
$ws.Cells.Item(167, 14).Value = 40000
$ws.Cells.Item(167, 15).Value = "00000000002083"
$ws.Cells.Item(167, 16).Value = "ALU-SERVIÇE"

# Row 168
$ws.Cells.Item(168, 1).Value = 2511
$ws.Cells.Item(168, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(168, 3).Value = "SP"
$ws.Cells.Item(168, 4).Value = 16
$ws.Cells.Item(168, 5).Value = 46014.6959672454
$ws.Cells.Item(168, 6).Value = 81396
$ws.Cells.Item(168, 7).Value = 46014
$ws.Cells.Item(168, 8).Value = "C.04.0100"
$ws.Cells.Item(168, 9).Value = "DESINFETANTE 5 L"
$ws.Cells.Item(168, 10).Value = "Apto"
$ws.Cells.Item(168, 11).Value = "UN"
$ws.Cells.Item(168, 12).Value = 2
$ws.Cells.Item(168, 13).Value = 13.5
$ws.Cells.Item(168, 14).Value = 27
$ws.Cells.Item(168, 15).Value = "00000000007786"
$ws.Cells.Item(168, 16).Value = "VILE EMBALAGENS"

# Row 169
$ws.Cells.Item(169, 1).Value = 2511
$ws.Cells.Item(169, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(169, 3).Value = "SP"
$ws.Cells.Item(169, 4).Value = 16
$ws.Cells.Item(169, 5).Value = 46014.6959672454
$ws.Cells.Item(169, 6).Value = 81396
$ws.Cells.Item(169, 7).Value = 46014
$ws.Cells.Item(169, 8).Value = "C.04.0023"
$ws.Cells.Item(169, 9).Value = "SACO PLÁSTICO P/ LIXO - 100 L C/ 100 UN REFORÇADO"
$ws.Cells.Item(169, 10).Value = "Apto"
$ws.Cells.Item(169, 11).Value = "UN"
$ws.Cells.Item(169, 12).Value = 2
$ws.Cells.Item(169, 13).Value = 68
$ws.Cells.Item(169, 14).Value = 136
$ws.Cells.Item(169, 15).Value = "00000000007786"
$ws.Cells.Item(169, 16).Value = "VILE EMBALAGENS"

# Row 170
$ws.Cells.Item(170, 1).Value = 2511
$ws.Cells.Item(170, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(170, 3).Value = "SP"
$ws.Cells.Item(170, 4).Value = 16
$ws.Cells.Item(170, 5).Value = 46014.6959672454
$ws.Cells.Item(170, 6).Value = 81397
$ws.Cells.Item(170, 7).Value = 46014
$ws.Cells.Item(170, 8).Value = "E.04.0580"
$ws.Cells.Item(170, 9).Value = "CAIXA PLASTICA PARA ARGAMASSA DE PVC 40 L"
$ws.Cells.Item(170, 10).Value = "Apto"
$ws.Cells.Item(170, 11).Value = "UN"
$ws.Cells.Item(170, 12).Value = 3
$ws.Cells.Item(170, 13).Value = 43.5
$ws.Cells.Item(170, 14).Value = 130.5
$ws.Cells.Item(170, 15).Value = "00000000008882"
$ws.Cells.Item(170, 16).Value = "GALPÃO DAS FERRAMENT"

# Row 171
$ws.Cells.Item(171, 1).Value = 2511
$ws.Cells.Item(171, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(171, 3).Value = "SP"
$ws.Cells.Item(171, 4).Value = 16
$ws.Cells.Item(171, 5).Value = 46014.6959672454
$ws.Cells.Item(171, 6).Value = 81397
$ws.Cells.Item(171, 7).Value = 46014
$ws.Cells.Item(171, 8).Value = "E.04.0011"
$ws.Cells.Item(171, 9).Value = "TALHADEIRA FORJADA EM AÇO REDONDO - 3/4 X 12''"
$ws.Cells.Item(171, 10).Value = "Apto"
$ws.Cells.Item(171, 11).Value = "UN"
$ws.Cells.Item(171, 12).Value = 2
$ws.Cells.Item(171, 13).Value = 10.5
$ws.Cells.Item(171, 14).Value = 21
$ws.Cells.Item(171, 15).Value = "00000000008882"
$ws.Cells.Item(171, 16).Value = "GALPÃO DAS FERRAMENT"

# Row 172
$ws.Cells.Item(172, 1).Value = 2511
$ws.Cells.Item(172, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(172, 3).Value = "SP"
$ws.Cells.Item(172, 4).Value = 16
$ws.Cells.Item(172, 5).Value = 46014.6959672454
$ws.Cells.Item(172, 6).Value = 81397
$ws.Cells.Item(172, 7).Value = 46014
$ws.Cells.Item(172, 8).Value = "E.04.0502"
$ws.Cells.Item(172, 9).Value = "CHAVE PHILLIPS - DIAM: 1/4'' - COMP: 200 MM"
$ws.Cells.Item(172, 10).Value = "Apto"
$ws.Cells.Item(172, 11).Value = "UN"
$ws.Cells.Item(172, 12).Value = 1
$ws.Cells.Item(172, 13).Value = 12.9
$ws.Cells.Item(172, 14).Value = 12.9
$ws.Cells.Item(172, 15).Value = "00000000008882"
$ws.Cells.Item(172, 16).Value = "GALPÃO DAS FERRAMENT"

# Row 173
$ws.Cells.Item(173, 1).Value = 2511
$ws.Cells.Item(173, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(173, 3).Value = "SP"
$ws.Cells.Item(173, 4).Value = 16
$ws.Cells.Item(173, 5).Value = 46014.6959672454
$ws.Cells.Item(173, 6).Value = 81397
$ws.Cells.Item(173, 7).Value = 46014
$ws.Cells.Item(173, 8).Value = "E.04.0028"
$ws.Cells.Item(173, 9).Value = "PROTEÇÃO PARA TALHADEIRA E PONTEIRO EMPUNHADURA"
$ws.Cells.Item(173, 10).Value = "Apto"
$ws.Cells.Item(173, 11).Value = "UN"
$ws.Cells.Item(173, 12).Value = 3
$ws.Cells.Item(173, 13).Value = 40.1
$ws.Cells.Item(173, 14).Value = 120.3
$ws.Cells.Item(173, 15).Value = "00000000008882"
$ws.Cells.Item(173, 16).Value = "GALPÃO DAS FERRAMENT"

# Row 174
$ws.Cells.Item(174, 1).Value = 2511
$ws.Cells.Item(174, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(174, 3).Value = "SP"
$ws.Cells.Item(174, 4).Value = 16
$ws.Cells.Item(174, 5).Value = 46014.6959672454
$ws.Cells.Item(174, 6).Value = 81397
$ws.Cells.Item(174, 7).Value = 46014
$ws.Cells.Item(174, 8).Value = "E.04.0469"
$ws.Cells.Item(174, 9).Value = "CHAVE DE FENDA - DIAM: 1/4'' - COMP: 200 MM"
$ws.Cells.Item(174, 10).Value = "Apto"
$ws.Cells.Item(174, 11).Value = "UN"
$ws.Cells.Item(174, 12).Value = 1
$ws.Cells.Item(174, 13).Value = 6.55
$ws.Cells.Item(174, 14).Value = 6.55
$ws.Cells.Item(174, 15).Value = "00000000008882"
$ws.Cells.Item(174, 16).Value = "GALPÃO DAS FERRAMENT"

# Row 175
$ws.Cells.Item(175, 1).Value = 2511
$ws.Cells.Item(175, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(175, 3).Value = "SP"
$ws.Cells.Item(175, 4).Value = 16
$ws.Cells.Item(175, 5).Value = 46014.6959672454
$ws.Cells.Item(175, 6).Value = 81398
$ws.Cells.Item(175, 7).Value = 46014
$ws.Cells.Item(175, 8).Value = "K.01.0266"
$ws.Cells.Item(175, 9).Value = "CABO FLEXÍVEL PP - 750 V - 2 X 2,5 MM²"
$ws.Cells.Item(175, 10).Value = "Apto"
$ws.Cells.Item(175, 11).Value = "M"
$ws.Cells.Item(175, 12).Value = 30
$ws.Cells.Item(175, 13).Value = 6.42
$ws.Cells.Item(175, 14).Value = 192.6
$ws.Cells.Item(175, 15).Value = "00000000008626"
$ws.Cells.Item(175, 16).Value = "INOVA MATERIAIS"

# Row 176
$ws.Cells.Item(176, 1).Value = 2511
$ws.Cells.Item(176, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(176, 3).Value = "SP"
$ws.Cells.Item(176, 4).Value = 16
$ws.Cells.Item(176, 5).Value = 46014.6959672454
$ws.Cells.Item(176, 6).Value = 81404
$ws.Cells.Item(176, 7).Value = 46014
$ws.Cells.Item(176, 8).Value = "O.01.0142"
$ws.Cells.Item(176, 9).Value = "TABUA DE PINUS  1"" X 12"""
$ws.Cells.Item(176, 10).Value = "Apto"
$ws.Cells.Item(176, 11).Value = "M"
$ws.Cells.Item(176, 12).Value = 27
$ws.Cells.Item(176, 13).Value = 14.9
$ws.Cells.Item(176, 14).Value = 402.3
$ws.Cells.Item(176, 15).Value = "00000000008356"
$ws.Cells.Item(176, 16).Value = "PALMAPLASTIC"

# Row 177
$ws.Cells.Item(177, 1).Value = 2511
$ws.Cells.Item(177, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(177, 3).Value = "SP"
$ws.Cells.Item(177, 4).Value = 16
$ws.Cells.Item(177, 5).Value = 46014.6959672454
$ws.Cells.Item(177, 6).Value = 81397
$ws.Cells.Item(177, 7).Value = 46014
$ws.Cells.Item(177, 8).Value = "R.02.0115"
$ws.Cells.Item(177, 9).Value = "TINTA SPRAY"
$ws.Cells.Item(177, 10).Value = "Apto"
$ws.Cells.Item(177, 11).Value = "UN"
$ws.Cells.Item(177, 12).Value = 10
$ws.Cells.Item(177, 13).Value = 20
$ws.Cells.Item(177, 14).Value = 200
$ws.Cells.Item(177, 15).Value = "00000000008882"
$ws.Cells.Item(177, 16).Value = "GALPÃO DAS FERRAMENT"

# Row 178
$ws.Cells.Item(178, 1).Value = 2511
$ws.Cells.Item(178, 2).Value = "1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(178, 3).Value = "SP"
$ws.Cells.Item(178, 4).Value = 16
$ws.Cells.Item(178, 5).Value = 46014.6959672454
$ws.Cells.Item(178, 6).Value = 81397
$ws.Cells.Item(178, 7).Value = 46014
$ws.Cells.Item(178, 8).Value = "W.01.0047"
$ws.Cells.Item(178, 9).Value = "PREGO DE AÇO COM CABEÇA 18 X 27"
$ws.Cells.Item(178, 10).Value = "Apto"
$ws.Cells.Item(178, 11).Value = "KG"
$ws.Cells.Item(178, 12).Value = 3
$ws.Cells.Item(178, 13).Value = 33
$ws.Cells.Item(178, 14).Value = 99
$ws.Cells.Item(178, 15).Value = "00000000008882"
$ws.Cells.Item(178, 16).Value = "GALPÃO DAS FERRAMENT"

# Restore column O's normal (General) style now that the text values are locked in,
# so it matches the formatting of every other FORNECEDOR_CDG cell in the sheet.
$ws.Range("O162").Copy()
$ws.Range("O163:O178").PasteSpecial(-4122)
